# "Generate Report for Handback" — record that handback has completed for
# both localized-language sheets (zh-cn, de-de) and refresh the status text
# on the Overview sheet.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9d45149c5e53e6e3b5023cf16ffa4834f4e2ac3c/e2e/"

$handedBackStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both language status columns move from "Ready for
# handoff" to "Handed back: in sync with en-US" for both files.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $handedBackStatus
$overview.Range("F2").Value = $handedBackStatus
$overview.Range("E3").Value = $handedBackStatus
$overview.Range("F3").Value = $handedBackStatus

# ---------------------------------------------------------------------
# Per-language detail sheets: zh-cn and de-de each report on the same
# two source files (rows 2 and 3). Handback is now complete: the
# status column updates, and the "Latest Target File" / "Latest
# Handback File" / "Latest Handback DateTime" columns (I/J/K) are
# populated for the first time.
# ---------------------------------------------------------------------
$languages = @(
    @{
        SheetName = "zh-cn"
        HandbackDateTime = "2016-08-19 14:53:52"
        Row2 = @{
            FileId = "15757890-fe2d-41c0-a634-369d4eb07159"
            TargetFile = "15757890-fe2d-41c0-a634-369d4eb07159.md"
            HandbackFile = "15757890-fe2d-41c0-a634-369d4eb07159.10510752c74be15c32553ce9fd1ac03717a7f737.zh-cn.xlf"
        }
        Row3 = @{
            FileId = "3df35001-bd95-4631-aa88-e5606593fdd5"
            TargetFile = "3df35001-bd95-4631-aa88-e5606593fdd5.md"
            HandbackFile = "3df35001-bd95-4631-aa88-e5606593fdd5.9a8a538af38b08302a0793ea752f30c3c80b4caf.zh-cn.xlf"
        }
    },
    @{
        SheetName = "de-de"
        HandbackDateTime = "2016-08-19 14:53:58"
        Row2 = @{
            FileId = "15757890-fe2d-41c0-a634-369d4eb07159"
            TargetFile = "15757890-fe2d-41c0-a634-369d4eb07159.md"
            HandbackFile = "15757890-fe2d-41c0-a634-369d4eb07159.10510752c74be15c32553ce9fd1ac03717a7f737.de-de.xlf"
        }
        Row3 = @{
            FileId = "3df35001-bd95-4631-aa88-e5606593fdd5"
            TargetFile = "3df35001-bd95-4631-aa88-e5606593fdd5.md"
            HandbackFile = "3df35001-bd95-4631-aa88-e5606593fdd5.9a8a538af38b08302a0793ea752f30c3c80b4caf.de-de.xlf"
        }
    }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.SheetName)

    foreach ($rowNum in 2, 3) {
        $rowInfo = if ($rowNum -eq 2) { $lang.Row2 } else { $lang.Row3 }

        # Status column (C)
        $ws.Range("C$rowNum").Value = $handedBackStatus

        # Latest Target File (I) — gets the same hyperlink treatment as
        # column A (linking back to the source .md file on GitHub).
        $iCell = $ws.Range("I$rowNum")
        $ws.Hyperlinks.Add($iCell, ($repoBase + $rowInfo.FileId + ".md"), $null, $null, $rowInfo.TargetFile) | Out-Null

        # Latest Handback File (J)
        $ws.Range("J$rowNum").Value = $rowInfo.HandbackFile

        # Latest Handback DateTime (K)
        $ws.Range("K$rowNum").Value = $lang.HandbackDateTime
    }
}

# ---------------------------------------------------------------------
# Widen the status / target-file / handback-file columns so the longer
# "Handed back: in sync with en-US" text and file names are fully
# visible.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527 - (5 / 6)
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527 - (5 / 6)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.SheetName)
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527 - (5 / 6)
    $ws.Columns.Item(9).ColumnWidth = 40 - (5 / 6)
    $ws.Columns.Item(10).ColumnWidth = 40 - (5 / 6)
}

Write-Output "Handback report generated"
